# Fruta / hortaliza, semanal
# Insert a new weekly record at row 165 (pushing the existing rows 165-227
# down to 166-228) for "Feria Lagunitas de Puerto Montt" - Mango.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data down by inserting a new row at 165.
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new weekly entry.
$ws.Cells.Item(165, 1).Value  = 4
$ws.Cells.Item(165, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(165, 3).Value  = "Los Lagos"
$ws.Cells.Item(165, 4).Value  = 44784
$ws.Cells.Item(165, 5).Value  = 10
$ws.Cells.Item(165, 6).Value  = "Fruta"
$ws.Cells.Item(165, 7).Value  = 100108
$ws.Cells.Item(165, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(165, 9).Value  = 100108002
$ws.Cells.Item(165, 10).Value = "Mango"
$ws.Cells.Item(165, 11).Value = "Sin especificar"
$ws.Cells.Item(165, 12).Value = "Primera"
$ws.Cells.Item(165, 13).Value = 140
$ws.Cells.Item(165, 14).Value = 13000
$ws.Cells.Item(165, 15).Value = 14000
$ws.Cells.Item(165, 16).Value = 13500
$ws.Cells.Item(165, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(165, 18).Value = "Brasil"
$ws.Cells.Item(165, 19).Value = 3375
$ws.Cells.Item(165, 20).Value = 4
